$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (clarify uniqueness / formatting requirements)
$ws.Range("H1").Value = "Processes(separate by comma)"

# Fill in the new Processes(separate by comma) column data for the two sample rows
$ws.Range("H3").Value = "cheatengine.exe"
$ws.Range("H2").Value = "mikesunique.exe"

$ws.Range("A1").Value = "Name(Must Be unique)"

# Widen column A so the longer header text fits
$ws.Columns.Item(1).ColumnWidth = 26.333333333333336

# Move/save the active selection to L8 (matches the saved view state)
$excel.Goto($ws.Range("L8"))
